# ReconnectEmailTemp.docx edits
$d = $word.ActiveDocument

function Replace-InRange($range, $old, $new) {
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# --- Paragraph 5: "Hope you had a wonderful summer! We over here on the ProspectAve Team have been working hard ..." ---
$p5 = $d.Paragraphs(5).Range
Replace-InRange $p5 " We over here on the ProspectAve Team have been working hard " "  Here on the ProspectAve Team we have been working hard "

# --- Paragraph 7: "Since going live, ProspectAve.io has accumulated over 6,800 page views and 2,300 unique users (almost half of the undergrad population)" ---
$p7 = $d.Paragraphs(7).Range
Replace-InRange $p7 "ProspectAve.io has accumulated" "ProspectAve has accumulated"
Replace-InRange $p7 "(almost half of the undergrad population)" "(nearly half of the undergrad population)"

# --- Paragraph 8: "Active monthly users peaked at 1,223 over LWNPARTIES weekend" ---
$p8 = $d.Paragraphs(8).Range
Replace-InRange $p8 "LWNPARTIES" "LAWNPARTIES"

# --- Paragraph 9: "25% of our users access the page through their phones" ---
$p9 = $d.Paragraphs(9).Range
Replace-InRange $p9 "25% of our users access the page through their phones" "25% of users utilize the ProspectAve mobile page"

# --- Paragraph 27 (duplicate "Since going live..." block, no bookmark/text beyond ProspectAve.io -> ProspectAve) ---
$p27 = $d.Paragraphs(27).Range
Replace-InRange $p27 "ProspectAve.io has accumulated" "ProspectAve has accumulated"

# This duplicate block never had the "6,800"/"2,300" figures bolded; the edit adds that emphasis.
$p27b = $d.Paragraphs(27).Range
if ($p27b.Find.Execute("6,800", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $p27b.Bold = 1
}
$p27c = $d.Paragraphs(27).Range
if ($p27c.Find.Execute("2,300", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $p27c.Bold = 1
}

# --- Paragraph 28: LWNPARTIES -> LAWNPARTIES ---
$p28 = $d.Paragraphs(28).Range
Replace-InRange $p28 "LWNPARTIES" "LAWNPARTIES"

# --- Paragraph 29: 25% of our users ... -> 25% of users utilize the ProspectAve mobile page ---
$p29 = $d.Paragraphs(29).Range
Replace-InRange $p29 "25% of our users access the page through their phones" "25% of users utilize the ProspectAve mobile page"

# --- Paragraph 13: big paragraph with several edits ---
$p13 = $d.Paragraphs(13).Range
Replace-InRange $p13 "With the new semester upcoming and frosh week around the corner, I know this is an important time of the year for" "With the upcoming new semester and Frosh Week around the corner, we know this is an important time of year for"
Replace-InRange $p13 "for your frosh week events" "for your Frosh Week events"
Replace-InRange $p13 "love to get you on a FaceTime/" "love to have a FaceTime/"
Replace-InRange $p13 "for like 5 min next week to make sure we are on the same page for frosh week." "for roughly 5 minutes next week to make sure we are on the same page for Frosh Week."

# The old zero-length _GoBack bookmark sat mid-paragraph (between "FaceTime/"
# and "Google Hangouts call"). Re-adding a bookmark with the same reserved
# name replaces the old one and moves it to the last-edited spot -- here,
# the already-empty paragraph 14 right after the big paragraph (matching
# Word's own behaviour of relocating _GoBack to the final edit position).
$p14 = $d.Paragraphs(14).Range
$d.Bookmarks.Add("_GoBack", $p14) | Out-Null

Write-Output "done"
